$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.403.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +12.54%  "
$ws.Range("D3").Value = "'1.826.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +9.35%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'230.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.17%  "
$ws.Range("D6").Value = "'0.574"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.81%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'31.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.38%  "
$ws.Range("D9").Value = "'46.61"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.63%  "
$ws.Range("D10").Value = "'0.289"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.85%  "
$ws.Range("D11").Value = "'0.0679"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.47%  "
$ws.Range("E12").Value = "  +3.17%  "
$ws.Range("D13").Value = "'2.088.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.31%  "
$ws.Range("D14").Value = "'1.838.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.97%  "
$ws.Range("E15").Value = "  +8.52%  "
$ws.Range("D16").Value = "'34.375.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +12.45%  "
$ws.Range("D17").Value = "'10.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.50%  "
$ws.Range("E18").Value = "  +7.68%  "
$ws.Range("D19").Value = "'70.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.06%  "
$ws.Range("D20").Value = "'258.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.72%  "
$ws.Range("D21").Value = "'0.0₃0759"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.63%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "'10.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  +2.96%  "
$ws.Range("E25").Value = "  +4.15%  "
$ws.Range("D26").Value = "'159.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").Value = "'16.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.46%  "
$ws.Range("E28").Value = "  +5.42%  "
$ws.Range("D29").Value = "'7.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.78%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").Value = "'3.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +13.29%  "
$ws.Range("D32").Value = "'0.0525"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.84%  "
$ws.Range("D33").Value = "'1.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.66%  "
$ws.Range("D34").Value = "'3.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.56%  "
$ws.Range("D35").Value = "'1.538.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.77%  "
$ws.Range("D36").Value = "'1.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.95%  "
$ws.Range("E37").Value = "  +6.35%  "
$ws.Range("E38").Value = "  +7.27%  "
$ws.Range("E39").Value = "  +8.03%  "
$ws.Range("D40").Value = "'84.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("E41").Value = "  +5.38%  "
$ws.Range("E42").Value = "  +2.86%  "
$ws.Range("D43").Value = "'0.916"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.72%  "
$ws.Range("E44").Value = "  +5.95%  "
$ws.Range("D45").Value = "'0.0529"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.27%  "
$ws.Range("E46").Value = "  +6.22%  "
$ws.Range("D47").Value = "'1.979.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.58%  "
$ws.Range("E48").Value = "  +5.95%  "
$ws.Range("D49").Value = "'12.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +18.69%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").Value = "'51.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.66%  "
